# Natmi following Dr Hou advice: expand Lrfn3-Lrfn3 LR-pair table to include the ECs cluster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("ECs", "Lrfn3", "Lrfn3", "ECs", 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 0.03752330582044444, 0.337709752384, 0.003184043235430514, 0.003184043235430514),
    @("ECs", "Lrfn3", "Lrfn3", "FAPs", 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 0.4241657791884444, 3.817492012696, 0.03599262246212408, 0.03599262246212408),
    @("ECs", "Lrfn3", "Lrfn3", "sCs", 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 0.2032955562515555, 1.829660006264, 0.01725066132960942, 0.01725066132960942),
    @("FAPs", "Lrfn3", "Lrfn3", "ECs", 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 0.4241657791884444, 3.817492012696, 0.03599262246212408, 0.03599262246212408),
    @("FAPs", "Lrfn3", "Lrfn3", "FAPs", 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 4.794796308605445, 43.15316677744899, 0.4068628394506831, 0.406862839450683),
    @("FAPs", "Lrfn3", "Lrfn3", "sCs", 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 2.298065592504555, 20.682590332541, 0.1950025469345249, 0.1950025469345248),
    @("sCs", "Lrfn3", "Lrfn3", "ECs", 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 2, 0.6666666666666666, 0.1937093333333333, 0.581128, 0.05642732702716401, 0.05642732702716401, 0.2032955562515555, 1.829660006264, 0.01725066132960942, 0.01725066132960942),
    @("sCs", "Lrfn3", "Lrfn3", "FAPs", 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 3, 1, 2.189702333333333, 6.569107, 0.637858008847332, 0.6378580088473319, 2.298065592504555, 20.682590332541, 0.1950025469345249, 0.1950025469345248),
    @("sCs", "Lrfn3", "Lrfn3", "sCs", 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 3, 1, 1.049487666666667, 3.148463, 0.305714664125504, 0.305714664125504, 1.101424362485444, 9.912819262369, 0.09346145586136974, 0.09346145586136974)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $rowData[$j]
    }
}
